$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("inputs")

# Step 1: delete row 17 (South Coast, Northern pike, Garibaldi Lake) - entirely removed
$ws.Rows.Item(17).Delete()

# Step 2: insert 2 blank rows at row 2 (two single inserts), pushing rows 2-16 down to 4-18
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

# Step 3: copy the "special" rows (now at 17 and 18) up into the new blank rows 2 and 3
for ($col = 2; $col -le 6; $col++) {
    $ws.Cells.Item(2, $col).Value2 = $ws.Cells.Item(17, $col).Value2
    $ws.Cells.Item(3, $col).Value2 = $ws.Cells.Item(18, $col).Value2
}

# Step 4: the new rows 2 and 3 picked up the bold header formatting via Insert(); reset to default.
# Column A in these rows should be entirely empty (no cell record at all), so Clear() it outright;
# columns B:E just need their formatting reset back to normal (values already set above).
$ws.Range("A2:A3").Clear()
$ws.Range("B2:E3").ClearFormats()

# Step 5: remove the F value (long explanatory text) from new row 2, but keep wrap-text formatting
$ws.Cells.Item(2, 6).Value2 = $null
$ws.Cells.Item(2, 6).WrapText = $true

# Step 6: row 2 should not retain the old custom row height (60) any more - default height
$ws.Rows.Item(2).RowHeight = $ws.Rows.Item(4).RowHeight

# Step 7: delete the now-duplicated old rows (17 and 18)
$ws.Rows.Item(18).Delete()
$ws.Rows.Item(17).Delete()

Write-Host "Final UsedRange: $($ws.UsedRange.Address())"
for ($r = 1; $r -le 16; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $c = $ws.Cells.Item($r, 3).Value2
    $d = $ws.Cells.Item($r, 4).Value2
    $e = $ws.Cells.Item($r, 5).Value2
    $f = $ws.Cells.Item($r, 6).Value2
    Write-Host "Row $r : A=$a B=$b C=$c D=$d E=$e F=$f"
}
